$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

$ws.Range("A32:D32").Copy()
$ws.Range("A33:D33").PasteSpecial(-4122)
$ws.Range("A32:D32").Copy()
$ws.Range("A34:D34").PasteSpecial(-4122)

$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = "testcase32_adls_parquet_to_delta_auto"
$ws.Cells.Item(33,4).Value = "Y"

$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = "testcase33_adls_parquet_adls_csv_auto"
$ws.Cells.Item(34,4).Value = "Y"

$ws.Range("C27:C34").Formula = '=_xlfn.CONCAT("test/testcases/",B27,".xlsx")'
